# DAC_listed_countries.xlsx hand edit
# The "upper_middle" column (D) had "Grenada" and "Iran" accidentally pasted
# together into a single cell ("Grenada Iran") at D40, instead of being two
# separate rows as in the source PDF. This splits them back into two cells,
# pushing every entry below down by one row (D41:D81 -> D42:D82).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift D41:D81 down to D42:D82, preserving values + formatting ---
# (Range.Copy into an overlapping/adjacent destination doesn't reliably
# clear already-blank source cells in the destination, so stage the move
# through a scratch column well off to the side of the used range.)
$ws.Range("D41:D81").Copy($ws.Range("Z41:Z81"))
$ws.Range("D41:D82").ClearContents()
$ws.Range("Z41:Z81").Copy($ws.Range("D42:D82"))
$ws.Range("Z41:Z81").EntireColumn.Delete()

# --- Split the merged text back into its two rows ---
$ws.Range("D40").Value = "Grenada"
$ws.Range("D41").Value = "Iran"

# --- Restore the cursor/selection to where the edit was made ---
$ws.Range("D42").Select()
